$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.341.58'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '1.786.46'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '226.12'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '32.59'
$ws.Range('E8').Value = '  +1.86%  '
$ws.Range('D9').Value = '0.294'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('D10').Value = '0.0687'
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').Value = '2.046.95'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.792.58'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '11.02'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').Value = '0.633'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('D16').Value = '34.343.87'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('D17').Value = '4.28'
$ws.Range('E17').Value = '  +2.42%  '
$ws.Range('D18').Value = '68.22'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0793'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '244.25'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').Value = '11.16'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '4.14'
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '2.07'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '165.99'
$ws.Range('E25').Value = '  +2.41%  '
$ws.Range('D26').Value = '7.29'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').Value = '16.49'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').Value = '0.115'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('E30').Value = '  +7.56%  '
$ws.Range('D31').Value = '0.0525'
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').Value = '3.80'
$ws.Range('E32').Value = '  +2.47%  '
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('D34').Value = '1.81'
$ws.Range('E34').Value = '  +0.91%  '
$ws.Range('D35').Value = '2.60'
$ws.Range('E35').Value = '  +5.03%  '
$ws.Range('D36').Value = '1.408.89'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('D37').Value = '0.679'
$ws.Range('E37').Value = '  +4.64%  '
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('D39').Value = '0.0190'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').Value = '83.86'
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').Value = '2.78'
$ws.Range('E42').Value = '  +2.53%  '
$ws.Range('D43').Value = '0.935'
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('D44').Value = '13.83'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('D45').Value = '0.0526'
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('E46').Value = '  +3.16%  '
$ws.Range('D47').Value = '6.05'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').Value = '1.947.56'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').Value = '105.05'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('E51').Value = '  -3.37%  '
